# 242: myr changes (#277)
# Add a new "Compliant" column (E) to the "Details" sheet of the
# model-year-report assessment template.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")

# New header cell E1 -> "Compliant" (adds a new shared string + extends
# the used range from A1:D1 to A1:E1).
$ws.Range("E1").Value = "Compliant"

# Match the column width that was also recorded for column F (17 chars)
# even though it holds no data yet.
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668
